# Auto-generated COM-interop script applying the LOQ4238.xlsx diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column definitions: split the redundant A:B column-width band
# so column B is governed solely by its own (60.71-wide) col entry,
# matching <col min="1" max="1".../><col min="2" max="2".../>.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth()

# --- Clear the region that gets restructured (old rows 13-23, plus the
# new row 24) so every cell below can be written from a clean slate.
$ws.Range("A13:C24").Clear()

function Set-Label($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $true
    $r.WrapText = $false
    $r.VerticalAlignment = -4160
}

function Set-Body($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $false
    $r.WrapText = $true
    $r.VerticalAlignment = -4160
}

function Set-Red($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $false
    $r.WrapText = $true
    $r.VerticalAlignment = -4160
    $r.Font.Color = 255
}

# --- Row 10 ---
Set-Label "A10" 'Objetivos:'
Set-Body "B10" 'Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, 
Aplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso
Desenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning).'
Set-Red "C10" 'Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, 
Aplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso
Desenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning).'
$ws.Rows.Item(10).RowHeight = 60

# --- Row 13 ---
Set-Body "B13" '5840560 - Marco Antonio Carvalho Pereira'
Set-Red "C13" '5840560 - Marco Antonio Carvalho Pereira'
$ws.Rows.Item(13).AutoFit()

# --- Row 14 ---
Set-Label "A14" 'Programa resumido:'
Set-Body "B14" 'Tópicos que abordem o tema do projeto de seu planejamento a execução.'
Set-Red "C14" 'Tópicos que abordem o tema do projeto de seu planejamento a execução.'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15 ---
Set-Label "A15" 'Short syllabus:'
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 ---
Set-Label "A16" 'Programa:'
Set-Body "B16" 'Assuntos Técnicos específicos relacionados com o tema do projeto.'
Set-Red "C16" 'Assuntos Técnicos específicos relacionados com o tema do projeto.'
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17 ---
Set-Label "A17" 'Syllabus:'
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18 ---
Set-Label "A18" 'Avaliação:'
$ws.Rows.Item(18).AutoFit()

# --- Row 19 ---
Set-Label "A19" 'Método:'
Set-Body "B19" 'O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.'
Set-Red "C19" 'O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.'
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20 ---
Set-Label "A20" 'Critério:'
Set-Body "B20" 'A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.'
Set-Red "C20" 'A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.'
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21 ---
Set-Label "A21" 'Norma de recuperação:'
Set-Body "B21" 'Não há recuperação'
Set-Red "C21" 'Não há recuperação'
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 ---
Set-Label "A22" 'Bibliografia:'
Set-Body "B22" 'Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.
Livros e Artigos científicos relacionados com o tema do projeto.'
Set-Red "C22" 'Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.
Livros e Artigos científicos relacionados com o tema do projeto.'
$ws.Rows.Item(22).RowHeight = 120

# --- Row 23 ---
Set-Label "A23" 'Requisitos:'
$ws.Rows.Item(23).AutoFit()

# --- Row 24 ---
Set-Body "B24" 'LOQ4237 -  Projeto Integrado de Engenharia de Produção II  (Requisito fraco)
'
Set-Red "C24" 'LOQ4237 -  Projeto Integrado de Engenharia de Produção II  (Requisito fraco)
'
$ws.Rows.Item(24).RowHeight = 30

